$wb = $excel.ActiveWorkbook

# --- PIR sheet: append rows 271-284 ---
$ws = $wb.Worksheets.Item("PIR")
$ws.Range("A271:A284").NumberFormat = "@"
$ws.Cells.Item(271, 1).Value = '2026-02-06'
$ws.Cells.Item(271, 2).Value = '10:02:38'
$ws.Cells.Item(271, 3).Value = '10:00'
$ws.Cells.Item(271, 4).Value = 'Bathroom'
$ws.Cells.Item(271, 5).Value = 'No Motion'
$ws.Cells.Item(271, 6).Value = 'Inactive'
$ws.Cells.Item(272, 1).Value = '2026-02-06'
$ws.Cells.Item(272, 2).Value = '10:02:40'
$ws.Cells.Item(272, 3).Value = '10:00'
$ws.Cells.Item(272, 4).Value = 'Bathroom'
$ws.Cells.Item(272, 5).Value = 'No Motion'
$ws.Cells.Item(272, 6).Value = 'Inactive'
$ws.Cells.Item(273, 1).Value = '2026-02-06'
$ws.Cells.Item(273, 2).Value = '10:02:42'
$ws.Cells.Item(273, 3).Value = '10:00'
$ws.Cells.Item(273, 4).Value = 'Bathroom'
$ws.Cells.Item(273, 5).Value = 'No Motion'
$ws.Cells.Item(273, 6).Value = 'Inactive'
$ws.Cells.Item(274, 1).Value = '2026-02-06'
$ws.Cells.Item(274, 2).Value = '10:02:47'
$ws.Cells.Item(274, 3).Value = '10:00'
$ws.Cells.Item(274, 4).Value = 'Bathroom'
$ws.Cells.Item(274, 5).Value = 'No Motion'
$ws.Cells.Item(274, 6).Value = 'Inactive'
$ws.Cells.Item(275, 1).Value = '2026-02-06'
$ws.Cells.Item(275, 2).Value = '10:02:52'
$ws.Cells.Item(275, 3).Value = '10:00'
$ws.Cells.Item(275, 4).Value = 'Bathroom'
$ws.Cells.Item(275, 5).Value = 'No Motion'
$ws.Cells.Item(275, 6).Value = 'Inactive'
$ws.Cells.Item(276, 1).Value = '2026-02-06'
$ws.Cells.Item(276, 2).Value = '10:02:57'
$ws.Cells.Item(276, 3).Value = '10:00'
$ws.Cells.Item(276, 4).Value = 'Bathroom'
$ws.Cells.Item(276, 5).Value = 'No Motion'
$ws.Cells.Item(276, 6).Value = 'Inactive'
$ws.Cells.Item(277, 1).Value = '2026-02-06'
$ws.Cells.Item(277, 2).Value = '10:03:02'
$ws.Cells.Item(277, 3).Value = '10:00'
$ws.Cells.Item(277, 4).Value = 'Bathroom'
$ws.Cells.Item(277, 5).Value = 'No Motion'
$ws.Cells.Item(277, 6).Value = 'Inactive'
$ws.Cells.Item(278, 1).Value = '2026-02-06'
$ws.Cells.Item(278, 2).Value = '10:03:08'
$ws.Cells.Item(278, 3).Value = '10:00'
$ws.Cells.Item(278, 4).Value = 'Bathroom'
$ws.Cells.Item(278, 5).Value = 'No Motion'
$ws.Cells.Item(278, 6).Value = 'Inactive'
$ws.Cells.Item(279, 1).Value = '2026-02-06'
$ws.Cells.Item(279, 2).Value = '10:03:13'
$ws.Cells.Item(279, 3).Value = '10:00'
$ws.Cells.Item(279, 4).Value = 'Bathroom'
$ws.Cells.Item(279, 5).Value = 'No Motion'
$ws.Cells.Item(279, 6).Value = 'Inactive'
$ws.Cells.Item(280, 1).Value = '2026-02-06'
$ws.Cells.Item(280, 2).Value = '10:03:18'
$ws.Cells.Item(280, 3).Value = '10:00'
$ws.Cells.Item(280, 4).Value = 'Bathroom'
$ws.Cells.Item(280, 5).Value = 'No Motion'
$ws.Cells.Item(280, 6).Value = 'Inactive'
$ws.Cells.Item(281, 1).Value = '2026-02-06'
$ws.Cells.Item(281, 2).Value = '10:03:23'
$ws.Cells.Item(281, 3).Value = '10:00'
$ws.Cells.Item(281, 4).Value = 'Bathroom'
$ws.Cells.Item(281, 5).Value = 'No Motion'
$ws.Cells.Item(281, 6).Value = 'Inactive'
$ws.Cells.Item(282, 1).Value = '2026-02-06'
$ws.Cells.Item(282, 2).Value = '10:03:28'
$ws.Cells.Item(282, 3).Value = '10:00'
$ws.Cells.Item(282, 4).Value = 'Bathroom'
$ws.Cells.Item(282, 5).Value = 'No Motion'
$ws.Cells.Item(282, 6).Value = 'Inactive'
$ws.Cells.Item(283, 1).Value = '2026-02-06'
$ws.Cells.Item(283, 2).Value = '10:03:33'
$ws.Cells.Item(283, 3).Value = '10:00'
$ws.Cells.Item(283, 4).Value = 'Bathroom'
$ws.Cells.Item(283, 5).Value = 'No Motion'
$ws.Cells.Item(283, 6).Value = 'Inactive'
$ws.Cells.Item(284, 1).Value = '2026-02-06'
$ws.Cells.Item(284, 2).Value = '10:03:38'
$ws.Cells.Item(284, 3).Value = '10:00'
$ws.Cells.Item(284, 4).Value = 'Bathroom'
$ws.Cells.Item(284, 5).Value = 'No Motion'
$ws.Cells.Item(284, 6).Value = 'Inactive'
$ws.Range("A271:A284").Style = "Normal"

# --- Humidity sheet: append rows 173-181 ---
$ws = $wb.Worksheets.Item("Humidity")
$ws.Range("A173:A181").NumberFormat = "@"
$ws.Range("E173:E181").NumberFormat = "@"
$ws.Cells.Item(173, 1).Value = '2026-02-06'
$ws.Cells.Item(173, 2).Value = '10:02:38'
$ws.Cells.Item(173, 3).Value = '10:00'
$ws.Cells.Item(173, 4).Value = 'Bathroom'
$ws.Cells.Item(173, 5).Value = '69.5%'
$ws.Cells.Item(173, 6).Value = 'Active'
$ws.Cells.Item(174, 1).Value = '2026-02-06'
$ws.Cells.Item(174, 2).Value = '10:02:40'
$ws.Cells.Item(174, 3).Value = '10:00'
$ws.Cells.Item(174, 4).Value = 'Bathroom'
$ws.Cells.Item(174, 5).Value = '68.6%'
$ws.Cells.Item(174, 6).Value = 'Active'
$ws.Cells.Item(175, 1).Value = '2026-02-06'
$ws.Cells.Item(175, 2).Value = '10:02:45'
$ws.Cells.Item(175, 3).Value = '10:00'
$ws.Cells.Item(175, 4).Value = 'Bathroom'
$ws.Cells.Item(175, 5).Value = '69.6%'
$ws.Cells.Item(175, 6).Value = 'Active'
$ws.Cells.Item(176, 1).Value = '2026-02-06'
$ws.Cells.Item(176, 2).Value = '10:02:55'
$ws.Cells.Item(176, 3).Value = '10:00'
$ws.Cells.Item(176, 4).Value = 'Bathroom'
$ws.Cells.Item(176, 5).Value = '69.7%'
$ws.Cells.Item(176, 6).Value = 'Active'
$ws.Cells.Item(177, 1).Value = '2026-02-06'
$ws.Cells.Item(177, 2).Value = '10:03:00'
$ws.Cells.Item(177, 3).Value = '10:00'
$ws.Cells.Item(177, 4).Value = 'Bathroom'
$ws.Cells.Item(177, 5).Value = '68.7%'
$ws.Cells.Item(177, 6).Value = 'Active'
$ws.Cells.Item(178, 1).Value = '2026-02-06'
$ws.Cells.Item(178, 2).Value = '10:03:05'
$ws.Cells.Item(178, 3).Value = '10:00'
$ws.Cells.Item(178, 4).Value = 'Bathroom'
$ws.Cells.Item(178, 5).Value = '69.7%'
$ws.Cells.Item(178, 6).Value = 'Active'
$ws.Cells.Item(179, 1).Value = '2026-02-06'
$ws.Cells.Item(179, 2).Value = '10:03:15'
$ws.Cells.Item(179, 3).Value = '10:00'
$ws.Cells.Item(179, 4).Value = 'Bathroom'
$ws.Cells.Item(179, 5).Value = '69.7%'
$ws.Cells.Item(179, 6).Value = 'Active'
$ws.Cells.Item(180, 1).Value = '2026-02-06'
$ws.Cells.Item(180, 2).Value = '10:03:20'
$ws.Cells.Item(180, 3).Value = '10:00'
$ws.Cells.Item(180, 4).Value = 'Bathroom'
$ws.Cells.Item(180, 5).Value = '68.7%'
$ws.Cells.Item(180, 6).Value = 'Active'
$ws.Cells.Item(181, 1).Value = '2026-02-06'
$ws.Cells.Item(181, 2).Value = '10:03:35'
$ws.Cells.Item(181, 3).Value = '10:00'
$ws.Cells.Item(181, 4).Value = 'Bathroom'
$ws.Cells.Item(181, 5).Value = '69.6%'
$ws.Cells.Item(181, 6).Value = 'Active'
$ws.Range("A173:A181").Style = "Normal"
$ws.Range("E173:E181").Style = "Normal"

# --- Temperature sheet: append rows 173-181 ---
$ws = $wb.Worksheets.Item("Temperature")
$ws.Range("A173:A181").NumberFormat = "@"
$ws.Cells.Item(173, 1).Value = '2026-02-06'
$ws.Cells.Item(173, 2).Value = '10:02:39'
$ws.Cells.Item(173, 3).Value = '10:00'
$ws.Cells.Item(173, 4).Value = 'Bathroom'
$ws.Cells.Item(173, 5).Value = '27.9C'
$ws.Cells.Item(173, 6).Value = 'Active'
$ws.Cells.Item(174, 1).Value = '2026-02-06'
$ws.Cells.Item(174, 2).Value = '10:02:41'
$ws.Cells.Item(174, 3).Value = '10:00'
$ws.Cells.Item(174, 4).Value = 'Bathroom'
$ws.Cells.Item(174, 5).Value = '27.8C'
$ws.Cells.Item(174, 6).Value = 'Active'
$ws.Cells.Item(175, 1).Value = '2026-02-06'
$ws.Cells.Item(175, 2).Value = '10:02:45'
$ws.Cells.Item(175, 3).Value = '10:00'
$ws.Cells.Item(175, 4).Value = 'Bathroom'
$ws.Cells.Item(175, 5).Value = '27.8C'
$ws.Cells.Item(175, 6).Value = 'Active'
$ws.Cells.Item(176, 1).Value = '2026-02-06'
$ws.Cells.Item(176, 2).Value = '10:02:55'
$ws.Cells.Item(176, 3).Value = '10:00'
$ws.Cells.Item(176, 4).Value = 'Bathroom'
$ws.Cells.Item(176, 5).Value = '27.8C'
$ws.Cells.Item(176, 6).Value = 'Active'
$ws.Cells.Item(177, 1).Value = '2026-02-06'
$ws.Cells.Item(177, 2).Value = '10:03:00'
$ws.Cells.Item(177, 3).Value = '10:00'
$ws.Cells.Item(177, 4).Value = 'Bathroom'
$ws.Cells.Item(177, 5).Value = '27.8C'
$ws.Cells.Item(177, 6).Value = 'Active'
$ws.Cells.Item(178, 1).Value = '2026-02-06'
$ws.Cells.Item(178, 2).Value = '10:03:05'
$ws.Cells.Item(178, 3).Value = '10:00'
$ws.Cells.Item(178, 4).Value = 'Bathroom'
$ws.Cells.Item(178, 5).Value = '27.8C'
$ws.Cells.Item(178, 6).Value = 'Active'
$ws.Cells.Item(179, 1).Value = '2026-02-06'
$ws.Cells.Item(179, 2).Value = '10:03:15'
$ws.Cells.Item(179, 3).Value = '10:00'
$ws.Cells.Item(179, 4).Value = 'Bathroom'
$ws.Cells.Item(179, 5).Value = '27.7C'
$ws.Cells.Item(179, 6).Value = 'Active'
$ws.Cells.Item(180, 1).Value = '2026-02-06'
$ws.Cells.Item(180, 2).Value = '10:03:20'
$ws.Cells.Item(180, 3).Value = '10:00'
$ws.Cells.Item(180, 4).Value = 'Bathroom'
$ws.Cells.Item(180, 5).Value = '27.8C'
$ws.Cells.Item(180, 6).Value = 'Active'
$ws.Cells.Item(181, 1).Value = '2026-02-06'
$ws.Cells.Item(181, 2).Value = '10:03:36'
$ws.Cells.Item(181, 3).Value = '10:00'
$ws.Cells.Item(181, 4).Value = 'Bathroom'
$ws.Cells.Item(181, 5).Value = '27.7C'
$ws.Cells.Item(181, 6).Value = 'Active'
$ws.Range("A173:A181").Style = "Normal"
